# Generate Report for Handback
# ------------------------------------------------------------
# This mirrors the "Generate Report for Handback" commit: the
# localization-status workbook is updated to reflect that the two
# tracked files have been handed back (in sync with en-US), with
# their "Latest Target File" / "Latest Handback File" columns filled
# in and the handback datetime recorded.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# ---- zh-cn sheet: fill in Latest Target File (F) / Latest Handback File (G) ----
# Row 2 -> 13993422-cdc3-4777-90f6-4d24486970a4
$zh.Range("F2").Value = "13993422-cdc3-4777-90f6-4d24486970a4.md"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7f48e026e71ee2cf1119d642953404bd7006d1/e2e/13993422-cdc3-4777-90f6-4d24486970a4.md", "", "", "13993422-cdc3-4777-90f6-4d24486970a4.md") | Out-Null

$zh.Range("G2").Value = "13993422-cdc3-4777-90f6-4d24486970a4.9a6cf952fa0db7c90379d4afac827bca37ad03ee.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7105b22fe1f0520a42e205b48c4e433e09c40d51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/13993422-cdc3-4777-90f6-4d24486970a4.9a6cf952fa0db7c90379d4afac827bca37ad03ee.zh-cn.xlf", "", "", "13993422-cdc3-4777-90f6-4d24486970a4.9a6cf952fa0db7c90379d4afac827bca37ad03ee.zh-cn.xlf") | Out-Null

$zh.Range("H2").Value = "2016-03-20 02:48:01"

# Row 3 -> 69659aef-8495-4dcc-997b-87449dc8a14c
$zh.Range("F3").Value = "69659aef-8495-4dcc-997b-87449dc8a14c.md"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7f48e026e71ee2cf1119d642953404bd7006d1/e2e/69659aef-8495-4dcc-997b-87449dc8a14c.md", "", "", "69659aef-8495-4dcc-997b-87449dc8a14c.md") | Out-Null

$zh.Range("G3").Value = "69659aef-8495-4dcc-997b-87449dc8a14c.32495087d5a8ef6627b09ecf1d45a027ab19d4bb.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7105b22fe1f0520a42e205b48c4e433e09c40d51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/69659aef-8495-4dcc-997b-87449dc8a14c.32495087d5a8ef6627b09ecf1d45a027ab19d4bb.zh-cn.xlf", "", "", "69659aef-8495-4dcc-997b-87449dc8a14c.32495087d5a8ef6627b09ecf1d45a027ab19d4bb.zh-cn.xlf") | Out-Null

$zh.Range("H3").Value = "2016-03-20 02:48:01"

# ---- de-de sheet: fill in Latest Target File (F) / Latest Handback File (G) ----
# Row 2 -> 13993422-cdc3-4777-90f6-4d24486970a4
$de.Range("F2").Value = "13993422-cdc3-4777-90f6-4d24486970a4.md"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7f48e026e71ee2cf1119d642953404bd7006d1/e2e/13993422-cdc3-4777-90f6-4d24486970a4.md", "", "", "13993422-cdc3-4777-90f6-4d24486970a4.md") | Out-Null

$de.Range("G2").Value = "13993422-cdc3-4777-90f6-4d24486970a4.9a6cf952fa0db7c90379d4afac827bca37ad03ee.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33805c3b0567d16a9f669d2a0562aefa00cd77d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/13993422-cdc3-4777-90f6-4d24486970a4.9a6cf952fa0db7c90379d4afac827bca37ad03ee.de-de.xlf", "", "", "13993422-cdc3-4777-90f6-4d24486970a4.9a6cf952fa0db7c90379d4afac827bca37ad03ee.de-de.xlf") | Out-Null

$de.Range("H2").Value = "2016-03-20 02:48:07"

# Row 3 -> 69659aef-8495-4dcc-997b-87449dc8a14c
$de.Range("F3").Value = "69659aef-8495-4dcc-997b-87449dc8a14c.md"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b7f48e026e71ee2cf1119d642953404bd7006d1/e2e/69659aef-8495-4dcc-997b-87449dc8a14c.md", "", "", "69659aef-8495-4dcc-997b-87449dc8a14c.md") | Out-Null

$de.Range("G3").Value = "69659aef-8495-4dcc-997b-87449dc8a14c.32495087d5a8ef6627b09ecf1d45a027ab19d4bb.de-de.xlf"
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33805c3b0567d16a9f669d2a0562aefa00cd77d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/69659aef-8495-4dcc-997b-87449dc8a14c.32495087d5a8ef6627b09ecf1d45a027ab19d4bb.de-de.xlf", "", "", "69659aef-8495-4dcc-997b-87449dc8a14c.32495087d5a8ef6627b09ecf1d45a027ab19d4bb.de-de.xlf") | Out-Null

$de.Range("H3").Value = "2016-03-20 02:48:07"

# ---- Match the look of the existing hyperlink columns (A, B, D) on the new F/G cells ----
$hyperFont = @($zh.Range("F2"), $zh.Range("G2"), $zh.Range("F3"), $zh.Range("G3"), `
               $de.Range("F2"), $de.Range("G2"), $de.Range("F3"), $de.Range("G3"))
foreach ($cell in $hyperFont) {
    $cell.Font.Underline = $true
    $cell.Font.Color = 15570276
}
